# Update stale Zoho-CRM test-data credentials across the workbook.
# The old test account (prajjawalmodi05@gmail.com / Password@1313) is
# replaced everywhere by the current one (prajjawalmodi94@gmail.com /
# Shinchain@1100); a couple of now-blank cells on the Login-Credentials
# sheet are normalised to a literal `""`; and the Delete-Lead sheet picks
# up a mailto hyperlink on the username cell and becomes the active tab.

$wb = $excel.ActiveWorkbook

# ---- Login-Credentials ----------------------------------------------
$ws1 = $wb.Worksheets.Item("Login-Credentials")
$ws1.Range("B2").Value = "prajjawalmodi94@gmail.com"
$ws1.Range("C3").Value = '""'
$ws1.Range("C4").Value = '""'
$ws1.Range("B5").Value = '""'
$ws1.Range("C5").Value = '""'
$ws1.Range("C2").Value = "Shinchain@1100"
$ws1.Range("B7").Value = "prajjawalmodi94@gmail.com"
$ws1.Range("C7").Value = "Shinchain@1100"

# ---- Create-Lead ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Create-Lead")
$ws2.Range("B2").Value = "prajjawalmodi94@gmail.com"
$ws2.Range("C2").Value = "Shinchain@1100"
$ws2.Range("B3").Value = "prajjawalmodi94@gmail.com"
$ws2.Range("C3").Value = "Shinchain@1100"
$ws2.Range("B4").Value = "prajjawalmodi94@gmail.com"
$ws2.Range("C4").Value = "Shinchain@1100"

# ---- Edit-Lead ----------------------------------------------------
$ws3 = $wb.Worksheets.Item("Edit-Lead")
$ws3.Range("B2").Value = "prajjawalmodi94@gmail.com"
$ws3.Range("C2").Value = "Shinchain@1100"
$ws3.Range("B3").Value = "prajjawalmodi94@gmail.com"
$ws3.Range("C3").Value = "Shinchain@1100"
$ws3.Range("B4").Value = "prajjawalmodi94@gmail.com"
$ws3.Range("C4").Value = "Shinchain@1100"
$ws3.Range("B5").Value = "prajjawalmodi94@gmail.com"
$ws3.Range("C5").Value = "Shinchain@1100"

# ---- Filter-Lead --------------------------------------------------
$ws4 = $wb.Worksheets.Item("Filter-Lead")
$ws4.Range("B2").Value = "prajjawalmodi94@gmail.com"
$ws4.Range("C2").Value = "Shinchain@1100"
$ws4.Range("B3").Value = "prajjawalmodi94@gmail.com"
$ws4.Range("C3").Value = "Shinchain@1100"

# ---- Delete-Lead --------------------------------------------------
$ws5 = $wb.Worksheets.Item("Delete-Lead")
$ws5.Range("B2").Value = "prajjawalmodi94@gmail.com"
$ws5.Range("C2").Value = "Shinchain@1100"
$ws5.Range("B3").Value = "prajjawalmodi94@gmail.com"
$ws5.Range("C3").Value = "Shinchain@1100"

# B2 becomes a live mailto link (matching the username now shown there),
# styled with the workbook's built-in Hyperlink cell style.
$ws5.Hyperlinks.Add($ws5.Range("B2"), "mailto:prajjawalmodi94@gmail.com")
$ws5.Range("B2").Style = "Hyperlink"

# This sheet was the one being edited, so it becomes the active tab/sheet.
$ws5.Activate()
$ws5.Range("C4").Select()

Write-Host "done"
